$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6205
$ws.Range("C21").Value = 979
$ws.Range("D21").Value = 5582916
$ws.Range("E21").Value = 899.7447219983884
$ws.Range("F21").Value = 7.706995313313669
$ws.Range("G21").Value = 3.270042194092837
$ws.Range("H21").Value = 27.43555539729538
